# Horarios Linea 141 - actualizacion de datos scrapeados (05:16:02)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (main schedule sheet): rows 8-41 get new values, rows 42-47
# are brand new rows appended at the bottom.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value2 = "Última actualización: 05:16:02"
$ws1.Range("A3").Value2 = "Total filas: 42"

$sheet1Data = @(
    @(8, "04:44:46", "04:46", "15_ABASTO", 2, "LP1912"),
    @(9, "04:44:46", "04:46", "215_EL PELIGRO", 2, "LP1912"),
    @(10, "03:52:04", "04:46", "215A_EL PATO", 54, "LP1912"),
    @(11, "04:32:18", "04:47", "215_EL PELIGRO", 15, "LP1912"),
    @(12, "04:52:24", "04:53", "11_ETCHEVERRY", 1, "LP1912"),
    @(13, "04:52:24", "04:54", "15_ABASTO", 2, "LP1912"),
    @(14, "04:13:31", "05:11", "17_ROMERO", 58, "LP1912"),
    @(15, "03:52:04", "05:16", "17_ROMERO", 84, "LP1912"),
    @(16, "04:52:24", "05:22", "23_HERNANDEZ", 30, "LP1912"),
    @(17, "05:16:02", "05:25", "23_HERNANDEZ", 9, "LP1912"),
    @(18, "04:44:46", "05:31", "81_EL PELIGRO", 47, "LP1912"),
    @(19, "05:16:02", "05:32", "81_EL PELIGRO", 16, "LP1912"),
    @(20, "03:52:04", "05:35", "215B_EL PATO", 103, "LP1912"),
    @(21, "05:16:02", "05:44", "14_ABASTO", 28, "LP1912"),
    @(22, "03:52:04", "05:46", "15_ABASTO", 114, "LP1912"),
    @(23, "04:32:18", "05:47", "14_ABASTO", 75, "LP1912"),
    @(24, "04:13:31", "05:50", "14_ABASTO", 97, "LP1912"),
    @(25, "04:44:46", "05:51", "17_ROMERO", 67, "LP1912"),
    @(26, "05:16:02", "05:52", "17_ROMERO", 36, "LP1912"),
    @(27, "04:44:46", "06:00", "16_SANTA ANA", 76, "LP1912"),
    @(28, "05:16:02", "06:01", "16_SANTA ANA", 45, "LP1912"),
    @(29, "04:44:46", "06:03", "10_OLMOS", 79, "LP1912"),
    @(30, "05:16:02", "06:04", "10_OLMOS", 48, "LP1912"),
    @(31, "04:44:46", "06:10", "215A_EL PATO", 86, "LP1912"),
    @(32, "05:16:02", "06:11", "215A_EL PATO", 55, "LP1912"),
    @(33, "04:32:18", "06:15", "17_ROMERO", 103, "LP1912"),
    @(34, "05:16:02", "06:24", "11_ETCHEVERRY", 68, "LP1912"),
    @(35, "05:16:02", "06:27", "23_HERNANDEZ", 71, "LP1912"),
    @(36, "04:44:46", "06:28", "17_ROMERO", 104, "LP1912"),
    @(37, "04:44:46", "06:30", "16_SANTA ANA", 106, "LP1912"),
    @(38, "05:16:02", "06:31", "16_SANTA ANA", 75, "LP1912"),
    @(39, "05:16:02", "06:31", "17X38_ROMERO", 75, "LP1912"),
    @(40, "04:52:24", "06:36", "17_ROMERO", 104, "LP1912"),
    @(41, "05:16:02", "06:39", "225_C ROCA-H SUR", 83, "LP1912"),
    @(42, "05:16:02", "06:50", "17_ROMERO", 94, "LP1912"),
    @(43, "05:16:02", "06:51", "215A_EL PATO", 95, "LP1912"),
    @(44, "05:16:02", "06:54", "14_ABASTO", 98, "LP1912"),
    @(45, "05:16:02", "07:04", "225_GOMEZ", 108, "LP1912"),
    @(46, "05:16:02", "07:07", "215C_EL PATO", 111, "LP1912"),
    @(47, "05:16:02", "07:14", "14X44_ABASTO", 118, "LP1912")
)

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value2 = $row[1]
    $ws1.Cells.Item($r, 2).Value2 = $row[2]
    $ws1.Cells.Item($r, 3).Value2 = $row[3]
    $ws1.Cells.Item($r, 4).Value2 = $row[4]
    $ws1.Cells.Item($r, 5).Value2 = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": rows 11-12 updated, row 13 is new.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value2 = "Última actualización: 05:16:02"
$ws2.Range("A3").Value2 = "Total filas: 8"

$sheet2Data = @(
    @(11, "05:16:02", "06:11", "215A_EL PATO", 55, "LP1912"),
    @(12, "05:16:02", "06:51", "215A_EL PATO", 95, "LP1912"),
    @(13, "05:16:02", "07:07", "215C_EL PATO", 111, "LP1912")
)

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value2 = $row[1]
    $ws2.Cells.Item($r, 2).Value2 = $row[2]
    $ws2.Cells.Item($r, 3).Value2 = $row[3]
    $ws2.Cells.Item($r, 4).Value2 = $row[4]
    $ws2.Cells.Item($r, 5).Value2 = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173": only the "Ultima actualizacion" timestamp changes.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value2 = "Última actualización: 05:16:02"
